$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "lal mirchi"
$ws.Range("B7").Value = "cupboard near chimney"

$ws.Range("A8").Value = "sugar"
$ws.Range("B8").Value = "cupboard near fridge"

$ws.Range("A9").Value = "dal"
$ws.Range("B9").Value = "cupboard over fridge"

$ws.Range("A10").Value = "cheeselings"
$ws.Range("B10").Value = "drawer near sliding drawer"

$ws.Range("A11").Value = "plates"
$ws.Range("B11").Value = "drawer under stove"

$ws.Range("A12").Value = "bournvita"
$ws.Range("B12").Value = "sliding drawer"

$ws.Range("A13").Select()
